# Auto-generated script applying 2022-04-19 daily data update
# to violent-crime-full-year.xlsx (column I = year 2022 running totals).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 1752
$ws.Range("H3").Value = 8348
$ws.Range("I3").Value = 1857
$ws.Range("G4").Value = 1423
$ws.Range("I4").Value = 460
$ws.Range("I5").Value = 166
$ws.Range("I6").Value = 2288
$ws.Range("G7").Value = 24642
$ws.Range("H7").Value = 25966
$ws.Range("I7").Value = 6523

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 25
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 20

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 86
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 41
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I4").Value = 31
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 231
$ws.Range("I8").Value = 402
$ws.Range("I10").Value = 52
$ws.Range("I11").Value = 114
$ws.Range("I18").Value = 51
$ws.Range("I19").Value = 188
$ws.Range("G20").Value = 619
$ws.Range("I20").Value = 178
$ws.Range("I29").Value = 418
$ws.Range("I30").Value = 20
$ws.Range("I31").Value = 65
$ws.Range("I33").Value = 298
$ws.Range("I36").Value = 84
$ws.Range("I42").Value = 219
$ws.Range("I43").Value = 59
$ws.Range("I47").Value = 48
$ws.Range("I49").Value = 36
$ws.Range("I51").Value = 64
$ws.Range("I52").Value = 134
$ws.Range("I54").Value = 143
$ws.Range("I55").Value = 75
$ws.Range("H63").Value = 183
$ws.Range("I65").Value = 149
$ws.Range("I67").Value = 251
$ws.Range("I76").Value = 106
$ws.Range("I78").Value = 85
$ws.Range("I83").Value = 124
$ws.Range("I86").Value = 40
$ws.Range("I90").Value = 76
$ws.Range("I91").Value = 76
$ws.Range("I96").Value = 86
$ws.Range("G101").Value = 24642
$ws.Range("H101").Value = 25966
$ws.Range("I101").Value = 6523

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 42
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 124

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 66
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 298

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I2").Value = 11
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 126
$ws.Range("I3").Value = 136
$ws.Range("I6").Value = 134
$ws.Range("I7").Value = 418

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 79
$ws.Range("I3").Value = 46
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 106

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I4").Value = 16
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 25
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I4").Value = 1
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("G4").Value = 25
$ws.Range("I6").Value = 67
$ws.Range("G7").Value = 619
$ws.Range("I7").Value = 178

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 22
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 14
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I6").Value = 138
$ws.Range("I7").Value = 402

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I2").Value = 4
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 24
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 59

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 68
$ws.Range("I7").Value = 231

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I4").Value = 1
$ws.Range("I7").Value = 31
